$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September entry was logged at the top of the list, so insert a new
# blank row before row 29 - this pushes the existing rows 29:44 down to 30:45
# (the last entry, "Broadband" in row 44, ends up in row 45).
$ws.Rows.Item(29).EntireRow.Insert()

# Populate the newly inserted row 29 with the latest September entry.
$ws.Range("R29").Value = "corporate internet share"
$ws.Range("S29").Value = "2024-09-03 19:22:58"
